# Calculando dinamicamente os dados da coluna 'Total' de cada aba da planilha
#
# For every monthly/report sheet in the workbook, the "Total" column (D) is
# recomputed as Quantidade (B) * Preço Unitário (C) instead of being a
# hard-coded literal. This turns each previously-static total into a live
# formula so it stays correct if the quantity or unit price changes.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Vendas", "Relatorio", "Janeiro", "Fevereiro", "Março")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    if ($lastRow -lt 2) {
        continue
    }

    for ($row = 2; $row -le $lastRow; $row++) {
        $ws.Cells.Item($row, 4).Formula = "=B" + $row + "*C" + $row
    }
}

Write-Host "Total column formulas applied to: $($sheetNames -join ', ')"
